# Updated symbol list on Thu Feb  2 11:09:59 UTC 2023 with GitHub Actions
# Refreshes Price (D), Volume(1h) (E) and Hora (G) columns for the crypto
# symbol table on the active sheet. Values are written with a leading
# apostrophe so Excel stores them as literal text (matching the sheet's
# existing text-formatted cells) instead of auto-converting to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ D = new-price; E = new-volume; G = new-hora }  (omitted keys are unchanged)
$updates = @{
    2 = @{ D="329.84"; E="7.31%"; G="11" }
    3 = @{ D="40.08"; E="8.25%"; G="11" }
    4 = @{ D="5.362"; E="4.79%"; G="11" }
    5 = @{ D="0.08104"; E="3.78%"; G="11" }
    6 = @{ D="4.528"; E="2.92%"; G="11" }
    7 = @{ D="8.649"; E="4.96%"; G="11" }
    8 = @{ D="1.921"; E="2.18%"; G="11" }
    9 = @{ G="11" }
    10 = @{ D="0.9431"; E="2.43%"; G="11" }
    11 = @{ D="0.1362"; E="25.98%"; G="11" }
    12 = @{ D="0.1968"; E="4.26%"; G="11" }
    13 = @{ D="0.09331"; E="5.12%"; G="11" }
    14 = @{ D="0.03565"; E="6.45%"; G="11" }
    15 = @{ D="0.09589"; E="0.17%"; G="11" }
    16 = @{ D="0.001317"; E="-4.62%"; G="11" }
    17 = @{ D="0.006352"; E="11.76%"; G="11" }
    18 = @{ D="3.369"; E="-1.30%"; G="11" }
    19 = @{ D="0.3521"; E="2.84%"; G="11" }
    20 = @{ D="7.194"; E="15.17%"; G="11" }
    21 = @{ D="0.1332"; E="2.13%"; G="11" }
    22 = @{ D="0.2561"; E="4.44%"; G="11" }
    23 = @{ D="0.04428"; E="2.23%"; G="11" }
    24 = @{ E="2.39%"; G="11" }
    25 = @{ D="0.004278"; E="0.42%"; G="11" }
    26 = @{ D="0.0001200"; E="-14.30%"; G="11" }
    27 = @{ D="0.0003991"; E="-0.02%"; G="11" }
    28 = @{ G="11" }
    29 = @{ G="11" }
    30 = @{ G="11" }
    31 = @{ G="11" }
    32 = @{ G="11" }
    33 = @{ G="11" }
    34 = @{ G="11" }
    35 = @{ G="11" }
    36 = @{ G="11" }
    37 = @{ G="11" }
    38 = @{ G="11" }
    39 = @{ D="0.02481"; E="15.01%"; G="11" }
    40 = @{ D="0.05230"; E="4.26%"; G="11" }
    41 = @{ D="0.007580"; E="0.90%"; G="11" }
    42 = @{ D="0.1427"; E="5.85%"; G="11" }
    43 = @{ D="0.009099"; E="5.13%"; G="11" }
    44 = @{ D="0.002171"; E="6.33%"; G="11" }
    45 = @{ D="0.01102"; E="26.42%"; G="11" }
    46 = @{ D="0.00006657"; E="1.86%"; G="11" }
    47 = @{ D="0.00000000750"; E="-0.01%"; G="11" }
    48 = @{ D="0.002401"; E="139.55%"; G="11" }
    49 = @{ E="1.39%"; G="11" }
    50 = @{ D="0.00002101"; E="-0.01%"; G="11" }
    51 = @{ D="0.0002001"; E="-0.01%"; G="11" }
}

foreach ($row in $updates.Keys) {
    $cells = $updates[$row]
    foreach ($col in $cells.Keys) {
        $ws.Range("$col$row").Value = "'" + $cells[$col]
    }
}

